$d = $word.ActiveDocument

# The paragraph currently reads "Version 2." built from the runs:
#   "Versi" | "on" | [spellEnd] | " 2" | [bookmarkStart _GoBack][bookmarkEnd] | "."
# It must become "Version 1." built from:
#   "Version" | [spellEnd] | " 1." | [bookmarkStart _GoBack][bookmarkEnd]
# (i.e. "Versi"+"on" merge into one run, " 2"+"." merge into one run and the
#  _GoBack bookmark ends up after all of the paragraph's text instead of
#  between the digit and the final period.)

# --- Step 0: merge "Versi" + "on" into a single "Version" run -------------
# A same-text assignment is a no-op for this engine, so first set it to a
# different value to force the run merge, then set it to the real value.
$r0 = $d.Range(0, 7)
$r0.Text = "VersionX"
$r0b = $d.Range(0, 8)
$r0b.Text = "Version"

# --- Step 1: change the digit "2" to "1" in place --------------------------
# This keeps the existing run/bookmark layout untouched (bookmark currently
# sits between " 1" and ".").
$rNum = $d.Range(8, 9)
$rNum.Text = "1"

# --- Step 2: temporarily insert a placeholder character after the "." -----
# Adding a zero-length bookmark exactly at the end of the paragraph's text
# (right before the paragraph mark) is unreliable, so we first push that
# position away from the paragraph edge with a throw-away character.
$insertPoint = $d.Range(10, 10)
$insertPoint.InsertAfter("Z")

# --- Step 3: (re)create the _GoBack bookmark at that now-safe position ----
$bmRange = $d.Range(10, 10)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Step 4: merge " 1" + "." into a single " 1." run ----------------------
# Same no-op-avoidance trick as step 0.
$r3 = $d.Range(7, 10)
$r3.Text = " 1.X"
$r4 = $d.Range(7, 11)
$r4.Text = " 1."

# --- Step 5: remove the placeholder character ------------------------------
$delRange = $d.Range(10, 11)
$delRange.Delete()
